$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "File Name" column header (H1) and its column width
$ws.Range("H1").Value = "File Name"
$ws.Columns.Item(8).ColumnWidth = 18

# New test-case row (row 4): Test Case id in A4, attached file name in H4
$ws.Range("A4").Value = "SU-T74"
$ws.Range("H4").Value = "village_details.pdf"

# Style A4 (Courier New, 10pt, green) to match the hyperlink-style look used
# for this test case id
$font = $ws.Range("A4").Font
$font.Name = "Courier New"
$font.Size = 10
$font.Color = 6258495
$font.Family = 3

$ws.Range("B7").Select()
